# hof_inducted.xlsx cleanup:
#  - the worksheet table ("hof_inducted") is no longer needed, so convert
#    it back into a plain range (drops xl/tables/table1.xml + its rels)
#  - rename the sheet tab from the default "Sheet1" to "hof_inducted"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove any Excel Table objects on the sheet, keeping their data/range in
# place (this is the equivalent of Table Tools > Convert to Range).
for ($i = $ws.ListObjects.Count; $i -ge 1; $i--) {
    $ws.ListObjects.Item($i).Unlist()
}

# Rename the worksheet to match the table/file name.
$ws.Name = "hof_inducted"
